$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)

# --- Row 1 (headers) ---
# The existing header cells B1:F1 mistakenly held copies of the data row; replace them
# with proper field names, matching the header convention used on the other sheets.
$ws3.Range("B1").Value = "bank"
$ws3.Range("C1").Value = "deposit_type"
$ws3.Range("D1").Value = "currency"
$ws3.Range("E1").Value = "owner"
$ws3.Range("F1").Value = "total"
$ws3.Range("G1").Value = "property_category"
$ws3.Range("H1").Value = "category"
$ws3.Range("I1").Value = "date"
$ws3.Range("J1").Value = "legislator_name"
$ws3.Range("K1").Value = "legislator_id"
$ws3.Range("L1").Value = "source_file"
$ws3.Range("M1").Value = "index"

# Copy header formatting (bold font + border) from B1 onto the newly used header cells
$ws3.Range("B1").Copy() | Out-Null
$ws3.Range("G1:M1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 2 (data) ---
# Force I2 to stay text (otherwise the date-looking string becomes a date serial number).
$ws3.Range("I2").NumberFormat = "@"

$ws3.Range("G2").Value = "deposit"
$ws3.Range("H2").Value = "normal"
$ws3.Range("I2").Value = "2012-04-30"
$ws3.Range("J2").Value = "高金素梅"
$ws3.Range("K2").Value = 926
$ws3.Range("L2").Value = "tmpb18e1"
$ws3.Range("M2").Value = 45

# Now copy the plain data formatting from B2 onto the new data cells, to match the
# rest of row 2 (plain, non-bold, no border, General number format).
$ws3.Range("B2").Copy() | Out-Null
$ws3.Range("G2:M2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0
